$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.866.50'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '2.545.96'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.05'
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.84'
$ws.Range("E6").Value = '  +6.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.578'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.88'
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0827'
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.63'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '2.934.70'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '2.529.75'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.08'
$ws.Range("E16").Value = '  +6.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.868'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '42.882.77'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.59'
$ws.Range("E19").Value = '  +5.05%  '
$ws.Range("D20").Value = '0.0₃0993'
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.61'
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.05'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.57'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.18'
$ws.Range("E26").Value = '  -3.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.03'
$ws.Range("E29").Value = '  +3.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.12'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  +3.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.13'
$ws.Range("E32").Value = '  +3.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.46'
$ws.Range("E33").Value = '  +14.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.14'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0801'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.32'
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("E37").Value = '  -4.49%  '
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.82'
$ws.Range("E39").Value = '  +6.15%  '
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.44'
$ws.Range("E41").Value = '  +1.55%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.08'
$ws.Range("E43").Value = '  +29.27%  '
$ws.Range("D44").Value = '2.099.63'
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.93'
$ws.Range("E47").Value = '  +5.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.97'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").Value = '2.792.12'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.91'
$ws.Range("E50").Value = '  +9.18%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.29'
$ws.Range("E51").Value = '  -0.69%  '
